$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC (50 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 252.11765
$ws.Range("I33").Value = 232.26666
$ws.Range("K33").Value = 232.26666
$ws.Range("M33").Value = -3.266660000000002
$ws.Range("H57").Value = 49980
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 49980
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 149940
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -150938
$ws.Range("H58").Value = 6101
$ws.Range("I58").Value = 601.44446
$ws.Range("J58").Value = 18475
$ws.Range("K58").Value = 1804.33338
$ws.Range("L58").Value = 55425
$ws.Range("M58").Value = -1654.33338
$ws.Range("N58").Value = -55725
$ws.Range("H132").Value = 24049616
$ws.Range("I132").Value = 29413292
$ws.Range("J132").Value = 1253988.4
$ws.Range("K132").Value = 88239876
$ws.Range("L132").Value = 3761965.2
$ws.Range("M132").Value = -88237346
$ws.Range("N132").Value = -3767025.2
$ws.Range("H135").Value = 657.1579
$ws.Range("I135").Value = 434.33334
$ws.Range("J135").Value = 1492.75
$ws.Range("K135").Value = 3909.00006
$ws.Range("L135").Value = 13434.75
$ws.Range("M135").Value = -1374.00006
$ws.Range("N135").Value = -18504.75
$ws.Range("H137").Value = 2804162
$ws.Range("I137").Value = 3970071.5
$ws.Range("J137").Value = 5979.6
$ws.Range("K137").Value = 11910214.5
$ws.Range("L137").Value = 17938.8
$ws.Range("M137").Value = -11907664.5
$ws.Range("N137").Value = -23038.8
$ws.Range("H138").Value = 2771.9773
$ws.Range("I138").Value = 1534.0667
$ws.Range("J138").Value = 3026.3425
$ws.Range("K138").Value = 4602.2001
$ws.Range("L138").Value = 9079.0275
$ws.Range("M138").Value = 537.7999
$ws.Range("N138").Value = -19359.0275
$ws.Range("H141").Value = 7056.8
$ws.Range("I141").Value = 7958.5
$ws.Range("K141").Value = 23875.5
$ws.Range("M141").Value = -18695.5

# --- Worksheet: ARM (12 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 34301
$ws.Range("J101").Value = 34301
$ws.Range("L101").Value = 34301
$ws.Range("N101").Value = -40791
$ws.Range("H133").Value = 24020
$ws.Range("J133").Value = 24020
$ws.Range("L133").Value = 24020
$ws.Range("N133").Value = -29080
$ws.Range("H137").Value = 42467.6
$ws.Range("J137").Value = 42467.6
$ws.Range("L137").Value = 42467.6
$ws.Range("N137").Value = -52667.6

# --- Worksheet: BSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 48928
$ws.Range("J137").Value = 48928
$ws.Range("L137").Value = 48928
$ws.Range("N137").Value = -59128
$ws.Range("H138").Value = 40862.4
$ws.Range("J138").Value = 40862.4
$ws.Range("L138").Value = 40862.4
$ws.Range("N138").Value = -51142.4

# --- Worksheet: CRP (59 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4380.8965
$ws.Range("I31").Value = 1532.8108
$ws.Range("J31").Value = 9398.951999999999
$ws.Range("K31").Value = 1532.8108
$ws.Range("L31").Value = 9398.951999999999
$ws.Range("M31").Value = -1237.8108
$ws.Range("N31").Value = -9988.951999999999
$ws.Range("H34").Value = 4380.8965
$ws.Range("I34").Value = 1532.8108
$ws.Range("J34").Value = 9398.951999999999
$ws.Range("K34").Value = 1532.8108
$ws.Range("L34").Value = 9398.951999999999
$ws.Range("M34").Value = -1330.8108
$ws.Range("N34").Value = -9802.951999999999
$ws.Range("H58").Value = 3184
$ws.Range("I58").Value = 1680
$ws.Range("J58").Value = 9200
$ws.Range("K58").Value = 1680
$ws.Range("L58").Value = 9200
$ws.Range("M58").Value = -1477
$ws.Range("N58").Value = -9606
$ws.Range("H98").Value = 47299
$ws.Range("J98").Value = 47299
$ws.Range("L98").Value = 47299
$ws.Range("N98").Value = -51791
$ws.Range("H134").Value = 4967.758
$ws.Range("I134").Value = 5219.077
$ws.Range("J134").Value = 4034.2856
$ws.Range("K134").Value = 15657.231
$ws.Range("L134").Value = 12102.8568
$ws.Range("M134").Value = -13122.231
$ws.Range("N134").Value = -17172.8568
$ws.Range("H136").Value = 3184
$ws.Range("I136").Value = 1680
$ws.Range("J136").Value = 9200
$ws.Range("K136").Value = 5040
$ws.Range("L136").Value = 27600
$ws.Range("M136").Value = -2490
$ws.Range("N136").Value = -32700
$ws.Range("H137").Value = 48613.332
$ws.Range("J137").Value = 48613.332
$ws.Range("L137").Value = 48613.332
$ws.Range("N137").Value = -58813.332
$ws.Range("H138").Value = 48822.5
$ws.Range("J138").Value = 48822.5
$ws.Range("L138").Value = 48822.5
$ws.Range("N138").Value = -59102.5
$ws.Range("H139").Value = 39000
$ws.Range("J139").Value = 39000
$ws.Range("L139").Value = 39000
$ws.Range("N139").Value = -49280
$ws.Range("H140").Value = 83725.38
$ws.Range("J140").Value = 89869.164
$ws.Range("L140").Value = 89869.164
$ws.Range("N140").Value = -100229.164
$ws.Range("H141").Value = 30664.285
$ws.Range("J141").Value = 30664.285
$ws.Range("L141").Value = 30664.285
$ws.Range("N141").Value = -41024.285

# --- Worksheet: CUL (14 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 160.25
$ws.Range("I33").Value = 160.25
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 961.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -678.5
$ws.Range("N33").ClearContents()
$ws.Range("H131").Value = 780
$ws.Range("I131").Value = 304.44446
$ws.Range("J131").Value = 827.03296
$ws.Range("K131").Value = 913.33338
$ws.Range("L131").Value = 2481.09888
$ws.Range("M131").Value = 4126.66662
$ws.Range("N131").Value = -12561.09888

# --- Worksheet: GSM (4 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 43780
$ws.Range("J137").Value = 43780
$ws.Range("L137").Value = 43780
$ws.Range("N137").Value = -53980

# --- Worksheet: LTW (84 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6668.1113
$ws.Range("J2").Value = 8001.7144
$ws.Range("L2").Value = 8001.7144
$ws.Range("N2").Value = -8225.714400000001
$ws.Range("H16").Value = 1111.5625
$ws.Range("I16").Value = 1105.6666
$ws.Range("K16").Value = 1105.6666
$ws.Range("M16").Value = -935.6666
$ws.Range("H22").Value = 35997.277
$ws.Range("I22").Value = 60244.707
$ws.Range("J22").Value = 1646.75
$ws.Range("K22").Value = 60244.707
$ws.Range("L22").Value = 1646.75
$ws.Range("M22").Value = -59949.707
$ws.Range("N22").Value = -2236.75
$ws.Range("H27").Value = 35997.277
$ws.Range("I27").Value = 60244.707
$ws.Range("J27").Value = 1646.75
$ws.Range("K27").Value = 60244.707
$ws.Range("L27").Value = 1646.75
$ws.Range("M27").Value = -60137.707
$ws.Range("N27").Value = -1860.75
$ws.Range("H40").Value = 9838.5
$ws.Range("I40").Value = 7740.6
$ws.Range("K40").Value = 7740.6
$ws.Range("M40").Value = -7604.6
$ws.Range("H46").Value = 1469.3846
$ws.Range("I46").Value = 1061.1111
$ws.Range("J46").Value = 2388
$ws.Range("K46").Value = 1061.1111
$ws.Range("L46").Value = 2388
$ws.Range("M46").Value = -873.1111000000001
$ws.Range("N46").Value = -2764
$ws.Range("H55").Value = 342.3846
$ws.Range("I55").Value = 107.6
$ws.Range("K55").Value = 107.6
$ws.Range("M55").Value = 65.40000000000001
$ws.Range("H61").Value = 1241.3704
$ws.Range("I61").Value = 1047.5294
$ws.Range("J61").Value = 1570.9
$ws.Range("K61").Value = 1047.5294
$ws.Range("L61").Value = 1570.9
$ws.Range("M61").Value = -845.5293999999999
$ws.Range("N61").Value = -1974.9
$ws.Range("H68").Value = 1182.6721
$ws.Range("I68").Value = 967.3273
$ws.Range("J68").Value = 3156.6667
$ws.Range("K68").Value = 967.3273
$ws.Range("L68").Value = 3156.6667
$ws.Range("M68").Value = -218.3273
$ws.Range("N68").Value = -4654.6667
$ws.Range("H71").Value = 1182.6721
$ws.Range("I71").Value = 967.3273
$ws.Range("J71").Value = 3156.6667
$ws.Range("K71").Value = 4836.636500000001
$ws.Range("L71").Value = 15783.3335
$ws.Range("M71").Value = -1092.636500000001
$ws.Range("N71").Value = -23271.3335
$ws.Range("H93").Value = 1716.4546
$ws.Range("I93").Value = 1192.4706
$ws.Range("J93").Value = 3498
$ws.Range("K93").Value = 1192.4706
$ws.Range("L93").Value = 3498
$ws.Range("M93").Value = 55.5293999999999
$ws.Range("N93").Value = -5994
$ws.Range("H113").Value = 1241.3704
$ws.Range("I113").Value = 1047.5294
$ws.Range("J113").Value = 1570.9
$ws.Range("K113").Value = 1047.5294
$ws.Range("L113").Value = 1570.9
$ws.Range("M113").Value = 1122.4706
$ws.Range("N113").Value = -5910.9
$ws.Range("H122").Value = 7831.1113
$ws.Range("I122").Value = 6782.857
$ws.Range("K122").Value = 20348.571
$ws.Range("M122").Value = -17898.571
$ws.Range("H139").Value = 47251
$ws.Range("J139").Value = 47251
$ws.Range("L139").Value = 47251
$ws.Range("N139").Value = -57531
$ws.Range("H141").Value = 32360
$ws.Range("J141").Value = 32360
$ws.Range("L141").Value = 32360
$ws.Range("N141").Value = -42720

# --- Worksheet: WVR (15 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2657.389
$ws.Range("I126").Value = 1131.8667
$ws.Range("J126").Value = 10285
$ws.Range("K126").Value = 3395.6001
$ws.Range("L126").Value = 30855
$ws.Range("M126").Value = -925.6001000000001
$ws.Range("N126").Value = -35795
$ws.Range("H133").Value = 43385
$ws.Range("J133").Value = 43385
$ws.Range("L133").Value = 43385
$ws.Range("N133").Value = -53505
$ws.Range("H136").Value = 3231.487
$ws.Range("J136").Value = 5822.6665
$ws.Range("L136").Value = 17467.9995
$ws.Range("N136").Value = -22567.9995
